# Generate Report for Handoff
# Update the "Latest Handoff Datetime" (column D) for the bc3b29b1-... source
# file (row 3) on both the zh-cn and de-de localization-status sheets, to
# reflect the timestamps recorded when the handoff report was generated.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-03-07 08:32:10"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-03-07 08:32:27"
